# Add a new worksheet "nr_studies" reporting the number of studies and
# effect sizes for each moderator level.

$wb = $excel.ActiveWorkbook

$lastIndex = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$ws.Name = "nr_studies"

# Header row.
$headers = @("outcome", "report_id_superior_report_e_g_parent_1_inferior_report_e_g_child_2_observation_3", "n_effect_sizes", "k_studies")
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}
$headerRange = $ws.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter

# Data rows. The moderator-level column (B) holds text digits ("1", "2", ...),
# so prefix with an apostrophe to force text entry rather than a number.
$data = @(
    @("NS", "1", 594, 71),
    @("NS", "2", 123, 20),
    @("NS", "3", 5, 1),
    @("NT", "1", 316, 44),
    @("NT", "2", 65, 12)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = "'" + $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}
